# Fix some dungeon config error (DungeonStory.xlsx)
# - I5: bossqiongqi=bossunicorn,manflower=wolfnest
#       -> bossqiongqi=bossunicorn,emanflower=fwolfnest
# - I6: bossqiongqi=forestexit,trees=rosemaryfield,river=poppyfield,manflower=trapspring,cliff=trappoison
#       -> bossqiongqi=forestexit,ftrees=rosemaryfield,river=poppyfield,emanflower=trapspring,cliff=trappoison
# - I8: bossmanwang=bosstalic,potteryroom=suntemple,trapspear=trapspearwall,stonedoor2=snowhill
#       -> bossmanwang=bosstalic,potteryroom=esuntemple,trapspear=trapspearwall,stonedoor2=snowhill
# - Window/background theme color (lt1) changed from white (FFFFFF) to C7EDCC
# - Active selection moved to M8

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("I5").Value = "bossqiongqi=bossunicorn,emanflower=fwolfnest"
$ws.Range("I6").Value = "bossqiongqi=forestexit,ftrees=rosemaryfield,river=poppyfield,emanflower=trapspring,cliff=trappoison"
$ws.Range("I8").Value = "bossmanwang=bosstalic,potteryroom=esuntemple,trapspear=trapspearwall,stonedoor2=snowhill"

# Theme "Background 1" (lt1) color -> C7EDCC (stored as BGR for the RGB() encoding)
$wb.Theme.ThemeColorScheme.Colors(2).RGB = 0xCCEDC7

# Move the active selection to M8, matching the saved view state
[void]$ws.Range("M8").Select()

Write-Host "edit applied"
